# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet right after "总计" (pushing 2022-Q3,
#    2022-Q2, 2022-Q1, 2021-Q4 down one slot each)
#  - add a new summary row for 2022-Q4 at the top of the "总计" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert the 2022-Q4 row, shift the rest down
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 6 is brand new (sheet used to stop at row 5) - give it the same
# style as row 5's A cell (style carrying the border/alignment/font) before
# filling it in, then write 2021-Q4's values (shifted down from row 5).
$total.Range("A6").Value = 4
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 0.1

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 25
$total.Range("D5").Value = 5.93

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 11
$total.Range("D4").Value = 2.98

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 34
$total.Range("D3").Value = 7.87

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 1.09

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the per-fund holdings table.
#    Clone the "2022-Q3" sheet (same column layout/styles) right after
#    "总计", trim it down to the rows we need, rename it, then overwrite
#    the data with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item("总计").Next()
$q4.Rows("9:35").Delete()
$q4.Name = "2022-Q4"

# Columns B-G hold text values (fund code / name / figures as strings);
# force text format before assigning so numeric-looking strings (e.g.
# "005106", "12.40") are not coerced into numbers and lose formatting.
$q4.Range("B2:G8").NumberFormat = "@"

$q4.Range("B2").Value = "005106"
$q4.Range("C2").Value = "银华农业产业股票A"
$q4.Range("D2").Value = "12.40"
$q4.Range("E2").Value = "93.03"
$q4.Range("F2").Value = "4.84"
$q4.Range("G2").Value = "0.6002"
$q4.Range("H2").Value = 9

$q4.Range("B3").Value = "164403"
$q4.Range("C3").Value = "前海开源沪港深农业混合（LOF）A"
$q4.Range("D3").Value = "4.16"
$q4.Range("E3").Value = "88.37"
$q4.Range("F3").Value = "6.48"
$q4.Range("G3").Value = "0.2696"
$q4.Range("H3").Value = 3

$q4.Range("B4").Value = "015210"
$q4.Range("C4").Value = "前海开源沪港深农业混合（LOF）C"
$q4.Range("D4").Value = "1.94"
$q4.Range("E4").Value = "88.37"
$q4.Range("F4").Value = "6.48"
$q4.Range("G4").Value = "0.1257"
$q4.Range("H4").Value = 3

$q4.Range("B5").Value = "014064"
$q4.Range("C5").Value = "银华农业产业股票C"
$q4.Range("D5").Value = "1.12"
$q4.Range("E5").Value = "93.03"
$q4.Range("F5").Value = "4.84"
$q4.Range("G5").Value = "0.0542"
$q4.Range("H5").Value = 9

$q4.Range("B6").Value = "004258"
$q4.Range("C6").Value = "国寿安保稳嘉混合A"
$q4.Range("D6").Value = "2.16"
$q4.Range("E6").Value = "23.32"
$q4.Range("F6").Value = "1.01"
$q4.Range("G6").Value = "0.0218"
$q4.Range("H6").Value = 10

$q4.Range("B7").Value = "006230"
$q4.Range("C7").Value = "鹏华研究驱动混合"
$q4.Range("D7").Value = "0.77"
$q4.Range("E7").Value = "87.24"
$q4.Range("F7").Value = "1.84"
$q4.Range("G7").Value = "0.0142"
$q4.Range("H7").Value = 10

$q4.Range("B8").Value = "004259"
$q4.Range("C8").Value = "国寿安保稳嘉混合C"
$q4.Range("D8").Value = "0.01"
$q4.Range("E8").Value = "23.32"
$q4.Range("F8").Value = "1.01"
$q4.Range("G8").Value = "0.0001"
$q4.Range("H8").Value = 10
